$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Drop the trailing listings (old rows 7-20). This also shrinks the
#    sheet's used range/dimension down to A1:H6 automatically.
# ---------------------------------------------------------------------
$ws.Rows("7:20").Delete()

# ---------------------------------------------------------------------
# 2) Refresh the remaining listing rows (2-6) with the new scrape data.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "2026-02-10 05:50:06"
$ws.Range("B2").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥API ◆ツール"

$ws.Range("A3").Value = "2026-02-10 05:50:06"
$ws.Range("B3").Value = "【急募】新聞記事PDFをCSV・Excel化するPythonプログラム作成依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5489128"
$ws.Range("G3").Value = 198
$ws.Range("H3").Value = "🔥Python"

$ws.Range("A4").Value = "2026-02-10 05:50:06"
$ws.Range("B4").Value = "【Java/講師】企業向け新入社員研修のJava講師業務(サブ講師)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5488955"
$ws.Range("G4").Value = 78
$ws.Range("H4").Value = "★Java"

$ws.Range("A5").Value = "2026-02-10 05:50:06"
$ws.Range("B5").Value = "スプレッドシート(Apps Script)で作業時間をボタン1つで計測・集計できる仕組みの開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5488743"
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = "◆開発"

$ws.Range("A6").Value = "2026-02-10 05:50:06"
$ws.Range("B6").Value = "【農機具管理】顧客指定で保有機情報を見れるシステム構築依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5489112"
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = "◇管理"

# ---------------------------------------------------------------------
# 3) Hyperlinks: drop every stale hyperlink (old F2..F20 set) then
#    re-create just the five that still matter (F2..F6), pointing to
#    the refreshed URLs above.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5489128")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5488955")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5488743")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5489112")

# ---------------------------------------------------------------------
# 4) Column width tweaks (B 52->49, D 30->28, H 19->12). ColumnWidth is
#    expressed in characters while the saved XML uses raw width units,
#    which run ~0.83 higher for this sheet's default font - compensate
#    so the stored width lands on the exact target value.
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 48.17
$ws.Columns("D").ColumnWidth = 27.17
$ws.Columns("H").ColumnWidth = 11.17
